$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to the "custom accuracy" (2-decimal) figures ---
$ws.Range("B5").Value = 22.24
$ws.Range("C5").Value = 16.8
$ws.Range("D5").Value = 0.92
$ws.Range("E5").Value = 47.9
$ws.Range("I5").Value = 26.51
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 18.12
$ws.Range("L5").Value = 19.37
$ws.Range("M5").Value = 20.34
$ws.Range("N5").Value = 5.53
$ws.Range("O5").Value = 17.18
$ws.Range("P5").Value = 24.69
$ws.Range("Q5").Value = 14.29
$ws.Range("R5").Value = 0.49
$ws.Range("S5").Value = 0.53
$ws.Range("T5").Value = 255.74
$ws.Range("U5").Value = 48.09
$ws.Range("V5").Value = 15.81
$ws.Range("W5").Value = 32.54
$ws.Range("X5").Value = 17.37
$ws.Range("Y5").Value = 2.27
$ws.Range("Z5").Value = 32.62
$ws.Range("AB5").Value = 12.87
$ws.Range("AC5").Value = 14.38
$ws.Range("AF5").Value = 58.47
$ws.Range("AG5").Value = 9.25
$ws.Range("AH5").Value = 19.8

# --- Delete row 6 (data trimmed to 1000 rows / 데이터 1000개) ---
$ws.Rows.Item(6).Delete()

# --- Narrow a set of columns from width 8 to width 7 ---
# (stored OOXML width = ColumnWidth + 0.8333333333333334, so target
#  stored width 7 => ColumnWidth 6.166666666666667)
$narrowWidth = 6.166666666666667
$narrowCols = @(3, 10, 11, 12, 13, 15, 17, 22, 24, 28, 29, 34)
foreach ($col in $narrowCols) {
    $ws.Columns.Item($col).ColumnWidth = $narrowWidth
}
